$wb = $excel.ActiveWorkbook

# --- Sheet1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M4").Value = 0
$ws1.Range("L8").Value = 0
$ws1.Range("D11").Value = 0
$ws1.Range("L18").Value = 0
$ws1.Range("L19").Value = 0
$ws1.Range("D24").Value = "0 de 22"
$ws1.Range("L24").Value = "0 de 22"
$ws1.Range("M24").Value = "0 de 22"

# --- Sheet2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
# Widen column F (6) from stored width 12 to 16 (ColumnWidth = width - 5/6)
$ws2.Columns.Item(6).ColumnWidth = 16 - (5/6)

# Month header shift
$ws2.Range("C1").Value = "junio"
$ws2.Range("D1").Value = "julio"
$ws2.Range("E1").Value = "agosto"
$ws2.Range("F1").Value = "septiembre"

# Data column shift (values move one month to the right)
$ws2.Range("D2").Value = 30.31
$ws2.Range("E2").Value = 0
$ws2.Range("D3").Value = 66.2
$ws2.Range("E3").Value = 0
$ws2.Range("E4").Value = 95.56
$ws2.Range("F4").Value = 0
$ws2.Range("C7").Value = 472.57
$ws2.Range("D7").Value = 0
$ws2.Range("E8").Value = 367.8
$ws2.Range("F8").Value = 0
$ws2.Range("D10").Value = 24.39
$ws2.Range("E10").Value = 0
$ws2.Range("E11").Value = 354.43
$ws2.Range("F11").Value = 0
$ws2.Range("C12").Value = 434.83
$ws2.Range("D12").Value = 0
$ws2.Range("D13").Value = 367.8
$ws2.Range("E13").Value = 0
$ws2.Range("C14").Value = 10.44
$ws2.Range("D14").Value = 0
$ws2.Range("C15").Value = 0
$ws2.Range("D16").Value = 238.35
$ws2.Range("E16").Value = 0
$ws2.Range("D17").Value = 308.08
$ws2.Range("E17").Value = 0
$ws2.Range("E18").Value = 238.35
$ws2.Range("F18").Value = 0
$ws2.Range("E19").Value = 367.8
$ws2.Range("F19").Value = 0
$ws2.Range("C20").Value = 159.03
$ws2.Range("D20").Value = 0
$ws2.Range("D21").Value = 413.5
$ws2.Range("E21").Value = 0
$ws2.Range("D22").Value = 43.86
$ws2.Range("E22").Value = 0
$ws2.Range("D23").Value = 565.27
$ws2.Range("E23").Value = 0
$ws2.Range("C24").Value = 1076.87
$ws2.Range("D24").Value = 2057.76
$ws2.Range("E24").Value = 1423.94
$ws2.Range("F24").Value = 0
